$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,7).Value = 0.7680745
$ws.Cells.Item(2,8).Value = 1.536149
$ws.Cells.Item(2,9).Value = 0.09258189785209336
$ws.Cells.Item(2,10).Value = 0.06368667635638427
$ws.Cells.Item(2,11).Value = 1.0
$ws.Cells.Item(2,12).Value = 0.5
$ws.Cells.Item(2,13).Value = 0.1443565
$ws.Cells.Item(2,14).Value = 0.288713
$ws.Cells.Item(2,15).Value = 0.001807064223911535
$ws.Cells.Item(2,16).Value = 0.001206425266199622
$ws.Cells.Item(2,17).Value = 0.11087654655925
$ws.Cells.Item(2,18).Value = 0.443506186237
$ws.Cells.Item(2,19).Value = 0.0001673014353903501
$ws.Cells.Item(2,20).Value = 0.00007683321547662005
$ws.Cells.Item(3,7).Value = 0.7680745
$ws.Cells.Item(3,8).Value = 1.536149
$ws.Cells.Item(3,9).Value = 0.09258189785209336
$ws.Cells.Item(3,10).Value = 0.06368667635638427
$ws.Cells.Item(3,15).Value = 0.0006061372836416816
$ws.Cells.Item(3,16).Value = 0.0006070005626485669
$ws.Cells.Item(3,17).Value = 0.0371909353645
$ws.Cells.Item(3,18).Value = 0.223145612187
$ws.Cells.Item(3,19).Value = 0.00005611734007845951
$ws.Cells.Item(3,20).Value = 0.00003865784838154243
$ws.Cells.Item(4,7).Value = 0.7680745
$ws.Cells.Item(4,8).Value = 1.536149
$ws.Cells.Item(4,9).Value = 0.09258189785209336
$ws.Cells.Item(4,10).Value = 0.06368667635638427
$ws.Cells.Item(4,13).Value = 27.803037
$ws.Cells.Item(4,14).Value = 83.409111
$ws.Cells.Item(4,15).Value = 0.3480402578255131
$ws.Cells.Item(4,16).Value = 0.3485359472612899
$ws.Cells.Item(4,17).Value = 21.3548037422565
$ws.Cells.Item(4,18).Value = 128.128822453539
$ws.Cells.Item(4,19).Value = 0.03222222759841789
$ws.Cells.Item(4,20).Value = 0.02219709607179559
$ws.Cells.Item(5,7).Value = 0.7680745
$ws.Cells.Item(5,8).Value = 1.536149
$ws.Cells.Item(5,9).Value = 0.09258189785209336
$ws.Cells.Item(5,10).Value = 0.06368667635638427
$ws.Cells.Item(5,13).Value = 0.19648
$ws.Cells.Item(5,14).Value = 0.39296
$ws.Cells.Item(5,15).Value = 0.002459549647671829
$ws.Cells.Item(5,16).Value = 0.001642035074990746
$ws.Cells.Item(5,17).Value = 0.15091127776
$ws.Cells.Item(5,18).Value = 0.60364511104
$ws.Cells.Item(5,19).Value = 0.0002277097742429054
$ws.Cells.Item(5,20).Value = 0.0001045757563867668
$ws.Cells.Item(6,7).Value = 0.7680745
$ws.Cells.Item(6,8).Value = 1.536149
$ws.Cells.Item(6,9).Value = 0.09258189785209336
$ws.Cells.Item(6,10).Value = 0.06368667635638427
$ws.Cells.Item(6,11).Value = 3.0
$ws.Cells.Item(6,12).Value = 1.0
$ws.Cells.Item(6,13).Value = 49.28281533333333
$ws.Cells.Item(6,14).Value = 147.848446
$ws.Cells.Item(6,15).Value = 0.6169255450395754
$ws.Cells.Item(6,16).Value = 0.6178041890138317
$ws.Cells.Item(6,17).Value = 37.85287374574233
$ws.Cells.Item(6,18).Value = 227.117242474454
$ws.Cells.Item(6,19).Value = 0.057116137793201
$ws.Cells.Item(6,20).Value = 0.03934589543734236
$ws.Cells.Item(7,7).Value = 0.7680745
$ws.Cells.Item(7,8).Value = 1.536149
$ws.Cells.Item(7,9).Value = 0.09258189785209336
$ws.Cells.Item(7,10).Value = 0.06368667635638427
$ws.Cells.Item(7,13).Value = 2.409433333333333
$ws.Cells.Item(7,14).Value = 7.2283
$ws.Cells.Item(7,15).Value = 0.03016144597968628
$ws.Cells.Item(7,16).Value = 0.03020440282103933
$ws.Cells.Item(7,17).Value = 1.850624302783333
$ws.Cells.Item(7,18).Value = 11.1037458167
$ws.Cells.Item(7,19).Value = 0.002792403910762747
$ws.Cells.Item(7,20).Value = 0.001923618027001391
$ws.Cells.Item(8,9).Value = 0.08079397830848588
$ws.Cells.Item(8,10).Value = 0.08336672828252452
$ws.Cells.Item(8,11).Value = 1.0
$ws.Cells.Item(8,12).Value = 0.5
$ws.Cells.Item(8,13).Value = 0.1443565
$ws.Cells.Item(8,14).Value = 0.288713
$ws.Cells.Item(8,15).Value = 0.001807064223911535
$ws.Cells.Item(8,16).Value = 0.001206425266199622
$ws.Cells.Item(8,17).Value = 0.09675927481999999
$ws.Cells.Item(8,18).Value = 0.58055564892
$ws.Cells.Item(8,19).Value = 0.0001459999077087494
$ws.Cells.Item(8,20).Value = 0.0001005757273604362
$ws.Cells.Item(9,9).Value = 0.08079397830848588
$ws.Cells.Item(9,10).Value = 0.08336672828252452
$ws.Cells.Item(9,15).Value = 0.0006061372836416816
$ws.Cells.Item(9,16).Value = 0.0006070005626485669
$ws.Cells.Item(9,19).Value = 0.00004897224254651058
$ws.Cells.Item(9,20).Value = 0.00005060365097366257
$ws.Cells.Item(10,9).Value = 0.08079397830848588
$ws.Cells.Item(10,10).Value = 0.08336672828252452
$ws.Cells.Item(10,13).Value = 27.803037
$ws.Cells.Item(10,14).Value = 83.409111
$ws.Cells.Item(10,15).Value = 0.3480402578255131
$ws.Cells.Item(10,16).Value = 0.3485359472612899
$ws.Cells.Item(10,17).Value = 18.63581964036
$ws.Cells.Item(10,18).Value = 167.72237676324
$ws.Cells.Item(10,19).Value = 0.02811955704123434
$ws.Cells.Item(10,20).Value = 0.02905630161202425
$ws.Cells.Item(11,9).Value = 0.08079397830848588
$ws.Cells.Item(11,10).Value = 0.08336672828252452
$ws.Cells.Item(11,13).Value = 0.19648
$ws.Cells.Item(11,14).Value = 0.39296
$ws.Cells.Item(11,15).Value = 0.002459549647671829
$ws.Cells.Item(11,16).Value = 0.001642035074990746
$ws.Cells.Item(11,17).Value = 0.1316966144
$ws.Cells.Item(11,18).Value = 0.7901796863999999
$ws.Cells.Item(11,19).Value = 0.0001987168008826418
$ws.Cells.Item(11,20).Value = 0.0001368910919271283
$ws.Cells.Item(12,9).Value = 0.08079397830848588
$ws.Cells.Item(12,10).Value = 0.08336672828252452
$ws.Cells.Item(12,11).Value = 3.0
$ws.Cells.Item(12,12).Value = 1.0
$ws.Cells.Item(12,13).Value = 49.28281533333333
$ws.Cells.Item(12,14).Value = 147.848446
$ws.Cells.Item(12,15).Value = 0.6169255450395754
$ws.Cells.Item(12,16).Value = 0.6178041890138317
$ws.Cells.Item(12,17).Value = 33.03328546162666
$ws.Cells.Item(12,18).Value = 297.29956915464
$ws.Cells.Item(12,19).Value = 0.04984386910387829
$ws.Cells.Item(12,20).Value = 0.05150431395732153
$ws.Cells.Item(13,9).Value = 0.08079397830848588
$ws.Cells.Item(13,10).Value = 0.08336672828252452
$ws.Cells.Item(13,13).Value = 2.409433333333333
$ws.Cells.Item(13,14).Value = 7.2283
$ws.Cells.Item(13,15).Value = 0.03016144597968628
$ws.Cells.Item(13,16).Value = 0.03020440282103933
$ws.Cells.Item(13,17).Value = 1.614994974666667
$ws.Cells.Item(13,18).Value = 14.534954772
$ws.Cells.Item(13,19).Value = 0.002436863212235341
$ws.Cells.Item(13,20).Value = 0.002518042242917502
$ws.Cells.Item(14,5).Value = 3.0
$ws.Cells.Item(14,6).Value = 1.0
$ws.Cells.Item(14,7).Value = 0.3750883333333333
$ws.Cells.Item(14,8).Value = 1.125265
$ws.Cells.Item(14,9).Value = 0.04521226751074096
$ws.Cells.Item(14,10).Value = 0.04665197703488838
$ws.Cells.Item(14,11).Value = 1.0
$ws.Cells.Item(14,12).Value = 0.5
$ws.Cells.Item(14,13).Value = 0.1443565
$ws.Cells.Item(14,14).Value = 0.288713
$ws.Cells.Item(14,15).Value = 0.001807064223911535
$ws.Cells.Item(14,16).Value = 0.001206425266199622
$ws.Cells.Item(14,17).Value = 0.05414643899083333
$ws.Cells.Item(14,18).Value = 0.324878633945
$ws.Cells.Item(14,19).Value = 0.00008170147110057781
$ws.Cells.Item(14,20).Value = 0.00005628212381305387
$ws.Cells.Item(15,5).Value = 3.0
$ws.Cells.Item(15,6).Value = 1.0
$ws.Cells.Item(15,7).Value = 0.3750883333333333
$ws.Cells.Item(15,8).Value = 1.125265
$ws.Cells.Item(15,9).Value = 0.04521226751074096
$ws.Cells.Item(15,10).Value = 0.04665197703488838
$ws.Cells.Item(15,15).Value = 0.0006061372836416816
$ws.Cells.Item(15,16).Value = 0.0006070005626485669
$ws.Cells.Item(15,17).Value = 0.01816215218833333
$ws.Cells.Item(15,18).Value = 0.163459369695
$ws.Cells.Item(15,19).Value = 0.00002740484101624158
$ws.Cells.Item(15,20).Value = 0.00002831777630884527
$ws.Cells.Item(16,5).Value = 3.0
$ws.Cells.Item(16,6).Value = 1.0
$ws.Cells.Item(16,7).Value = 0.3750883333333333
$ws.Cells.Item(16,8).Value = 1.125265
$ws.Cells.Item(16,9).Value = 0.04521226751074096
$ws.Cells.Item(16,10).Value = 0.04665197703488838
$ws.Cells.Item(16,13).Value = 27.803037
$ws.Cells.Item(16,14).Value = 83.409111
$ws.Cells.Item(16,15).Value = 0.3480402578255131
$ws.Cells.Item(16,16).Value = 0.3485359472612899
$ws.Cells.Item(16,17).Value = 10.428594809935
$ws.Cells.Item(16,18).Value = 93.85735328941499
$ws.Cells.Item(16,19).Value = 0.01573568924131435
$ws.Cells.Item(16,20).Value = 0.01625989100746677
$ws.Cells.Item(17,5).Value = 3.0
$ws.Cells.Item(17,6).Value = 1.0
$ws.Cells.Item(17,7).Value = 0.3750883333333333
$ws.Cells.Item(17,8).Value = 1.125265
$ws.Cells.Item(17,9).Value = 0.04521226751074096
$ws.Cells.Item(17,10).Value = 0.04665197703488838
$ws.Cells.Item(17,13).Value = 0.19648
$ws.Cells.Item(17,14).Value = 0.39296
$ws.Cells.Item(17,15).Value = 0.002459549647671829
$ws.Cells.Item(17,16).Value = 0.001642035074990746
$ws.Cells.Item(17,17).Value = 0.07369735573333332
$ws.Cells.Item(17,18).Value = 0.4421841344
$ws.Cells.Item(17,19).Value = 0.0001112018166264874
$ws.Cells.Item(17,20).Value = 0.00007660418260894953
$ws.Cells.Item(18,5).Value = 3.0
$ws.Cells.Item(18,6).Value = 1.0
$ws.Cells.Item(18,7).Value = 0.3750883333333333
$ws.Cells.Item(18,8).Value = 1.125265
$ws.Cells.Item(18,9).Value = 0.04521226751074096
$ws.Cells.Item(18,10).Value = 0.04665197703488838
$ws.Cells.Item(18,11).Value = 3.0
$ws.Cells.Item(18,12).Value = 1.0
$ws.Cells.Item(18,13).Value = 49.28281533333333
$ws.Cells.Item(18,14).Value = 147.848446
$ws.Cells.Item(18,15).Value = 0.6169255450395754
$ws.Cells.Item(18,16).Value = 0.6178041890138317
$ws.Cells.Item(18,17).Value = 18.48540906535444
$ws.Cells.Item(18,18).Value = 166.36868158819
$ws.Cells.Item(18,19).Value = 0.02789260277653896
$ws.Cells.Item(18,20).Value = 0.02882178683793112
$ws.Cells.Item(19,5).Value = 3.0
$ws.Cells.Item(19,6).Value = 1.0
$ws.Cells.Item(19,7).Value = 0.3750883333333333
$ws.Cells.Item(19,8).Value = 1.125265
$ws.Cells.Item(19,9).Value = 0.04521226751074096
$ws.Cells.Item(19,10).Value = 0.04665197703488838
$ws.Cells.Item(19,13).Value = 2.409433333333333
$ws.Cells.Item(19,14).Value = 7.2283
$ws.Cells.Item(19,15).Value = 0.03016144597968628
$ws.Cells.Item(19,16).Value = 0.03020440282103933
$ws.Cells.Item(19,17).Value = 0.9037503332777777
$ws.Cells.Item(19,18).Value = 8.1337529995
$ws.Cells.Item(19,19).Value = 0.001363667364144338
$ws.Cells.Item(19,20).Value = 0.001409095106759644
$ws.Cells.Item(20,7).Value = 6.453176333333333
$ws.Cells.Item(20,8).Value = 19.359529
$ws.Cells.Item(20,9).Value = 0.7778507320763975
$ws.Cells.Item(20,10).Value = 0.8026200959900606
$ws.Cells.Item(20,11).Value = 1.0
$ws.Cells.Item(20,12).Value = 0.5
$ws.Cells.Item(20,13).Value = 0.1443565
$ws.Cells.Item(20,14).Value = 0.288713
$ws.Cells.Item(20,15).Value = 0.001807064223911535
$ws.Cells.Item(20,16).Value = 0.001206425266199622
$ws.Cells.Item(20,17).Value = 0.9315579493628332
$ws.Cells.Item(20,18).Value = 5.589347696177
$ws.Cells.Item(20,19).Value = 0.001405626229478654
$ws.Cells.Item(20,20).Value = 0.0009683011629619749
$ws.Cells.Item(21,7).Value = 6.453176333333333
$ws.Cells.Item(21,8).Value = 19.359529
$ws.Cells.Item(21,9).Value = 0.7778507320763975
$ws.Cells.Item(21,10).Value = 0.8026200959900606
$ws.Cells.Item(21,15).Value = 0.0006061372836416816
$ws.Cells.Item(21,16).Value = 0.0006070005626485669
$ws.Cells.Item(21,17).Value = 0.3124692512363333
$ws.Cells.Item(21,18).Value = 2.812223261127
$ws.Cells.Item(21,19).Value = 0.0004714843298194811
$ws.Cells.Item(21,20).Value = 0.0004871908498590135
$ws.Cells.Item(22,7).Value = 6.453176333333333
$ws.Cells.Item(22,8).Value = 19.359529
$ws.Cells.Item(22,9).Value = 0.7778507320763975
$ws.Cells.Item(22,10).Value = 0.8026200959900606
$ws.Cells.Item(22,13).Value = 27.803037
$ws.Cells.Item(22,14).Value = 83.409111
$ws.Cells.Item(22,15).Value = 0.3480402578255131
$ws.Cells.Item(22,16).Value = 0.3485359472612899
$ws.Cells.Item(22,17).Value = 179.417900363191
$ws.Cells.Item(22,18).Value = 1614.761103268719
$ws.Cells.Item(22,19).Value = 0.2707233693416335
$ws.Cells.Item(22,20).Value = 0.2797419554468432
$ws.Cells.Item(23,7).Value = 6.453176333333333
$ws.Cells.Item(23,8).Value = 19.359529
$ws.Cells.Item(23,9).Value = 0.7778507320763975
$ws.Cells.Item(23,10).Value = 0.8026200959900606
$ws.Cells.Item(23,13).Value = 0.19648
$ws.Cells.Item(23,14).Value = 0.39296
$ws.Cells.Item(23,15).Value = 0.002459549647671829
$ws.Cells.Item(23,16).Value = 0.001642035074990746
$ws.Cells.Item(23,17).Value = 1.267920085973333
$ws.Cells.Item(23,18).Value = 7.607520515839999
$ws.Cells.Item(23,19).Value = 0.001913162494019777
$ws.Cells.Item(23,20).Value = 0.001317930349508119
$ws.Cells.Item(24,7).Value = 6.453176333333333
$ws.Cells.Item(24,8).Value = 19.359529
$ws.Cells.Item(24,9).Value = 0.7778507320763975
$ws.Cells.Item(24,10).Value = 0.8026200959900606
$ws.Cells.Item(24,11).Value = 3.0
$ws.Cells.Item(24,12).Value = 1.0
$ws.Cells.Item(24,13).Value = 49.28281533333333
$ws.Cells.Item(24,14).Value = 147.848446
$ws.Cells.Item(24,15).Value = 0.6169255450395754
$ws.Cells.Item(24,16).Value = 0.6178041890138317
$ws.Cells.Item(24,17).Value = 318.0306975491037
$ws.Cells.Item(24,18).Value = 2862.276277941934
$ws.Cells.Item(24,19).Value = 0.4798759868456643
$ws.Cells.Item(24,20).Value = 0.4958620574893431
$ws.Cells.Item(25,7).Value = 6.453176333333333
$ws.Cells.Item(25,8).Value = 19.359529
$ws.Cells.Item(25,9).Value = 0.7778507320763975
$ws.Cells.Item(25,10).Value = 0.8026200959900606
$ws.Cells.Item(25,13).Value = 2.409433333333333
$ws.Cells.Item(25,14).Value = 7.2283
$ws.Cells.Item(25,15).Value = 0.03016144597968628
$ws.Cells.Item(25,16).Value = 0.03020440282103933
$ws.Cells.Item(25,17).Value = 15.54849816341111
$ws.Cells.Item(25,18).Value = 139.9364834707
$ws.Cells.Item(25,19).Value = 0.02346110283578169
$ws.Cells.Item(25,20).Value = 0.02424266069154504
$ws.Cells.Item(26,5).Value = 1.0
$ws.Cells.Item(26,6).Value = 0.3333333333333333
$ws.Cells.Item(26,7).Value = 0.02954366666666667
$ws.Cells.Item(26,8).Value = 0.088631
$ws.Cells.Item(26,9).Value = 0.003561124252282336
$ws.Cells.Item(26,10).Value = 0.003674522336142324
$ws.Cells.Item(26,11).Value = 1.0
$ws.Cells.Item(26,12).Value = 0.5
$ws.Cells.Item(26,13).Value = 0.1443565
$ws.Cells.Item(26,14).Value = 0.288713
$ws.Cells.Item(26,15).Value = 0.001807064223911535
$ws.Cells.Item(26,16).Value = 0.001206425266199622
$ws.Cells.Item(26,17).Value = 0.004264820317166667
$ws.Cells.Item(26,18).Value = 0.025588921903
$ws.Cells.Item(26,19).Value = 0.000006435180233203123
$ws.Cells.Item(26,20).Value = 0.00000443303658753696
$ws.Cells.Item(27,5).Value = 1.0
$ws.Cells.Item(27,6).Value = 0.3333333333333333
$ws.Cells.Item(27,7).Value = 0.02954366666666667
$ws.Cells.Item(27,8).Value = 0.088631
$ws.Cells.Item(27,9).Value = 0.003561124252282336
$ws.Cells.Item(27,10).Value = 0.003674522336142324
$ws.Cells.Item(27,15).Value = 0.0006061372836416816
$ws.Cells.Item(27,16).Value = 0.0006070005626485669
$ws.Cells.Item(27,17).Value = 0.001430533883666667
$ws.Cells.Item(27,18).Value = 0.012874804953
$ws.Cells.Item(27,19).Value = 0.00000215853018098893
$ws.Cells.Item(27,20).Value = 0.000002230437125503117
$ws.Cells.Item(28,5).Value = 1.0
$ws.Cells.Item(28,6).Value = 0.3333333333333333
$ws.Cells.Item(28,7).Value = 0.02954366666666667
$ws.Cells.Item(28,8).Value = 0.088631
$ws.Cells.Item(28,9).Value = 0.003561124252282336
$ws.Cells.Item(28,10).Value = 0.003674522336142324
$ws.Cells.Item(28,13).Value = 27.803037
$ws.Cells.Item(28,14).Value = 83.409111
$ws.Cells.Item(28,15).Value = 0.3480402578255131
$ws.Cells.Item(28,16).Value = 0.3485359472612899
$ws.Cells.Item(28,17).Value = 0.821403657449
$ws.Cells.Item(28,18).Value = 7.392632917040999
$ws.Cells.Item(28,19).Value = 0.001239414602913032
$ws.Cells.Item(28,20).Value = 0.001280703123160133
$ws.Cells.Item(29,5).Value = 1.0
$ws.Cells.Item(29,6).Value = 0.3333333333333333
$ws.Cells.Item(29,7).Value = 0.02954366666666667
$ws.Cells.Item(29,8).Value = 0.088631
$ws.Cells.Item(29,9).Value = 0.003561124252282336
$ws.Cells.Item(29,10).Value = 0.003674522336142324
$ws.Cells.Item(29,13).Value = 0.19648
$ws.Cells.Item(29,14).Value = 0.39296
$ws.Cells.Item(29,15).Value = 0.002459549647671829
$ws.Cells.Item(29,16).Value = 0.001642035074990746
$ws.Cells.Item(29,17).Value = 0.005804739626666666
$ws.Cells.Item(29,18).Value = 0.03482843776
$ws.Cells.Item(29,19).Value = 0.000008758761900016623
$ws.Cells.Item(29,20).Value = 0.000006033694559782634
$ws.Cells.Item(30,5).Value = 1.0
$ws.Cells.Item(30,6).Value = 0.3333333333333333
$ws.Cells.Item(30,7).Value = 0.02954366666666667
$ws.Cells.Item(30,8).Value = 0.088631
$ws.Cells.Item(30,9).Value = 0.003561124252282336
$ws.Cells.Item(30,10).Value = 0.003674522336142324
$ws.Cells.Item(30,11).Value = 3.0
$ws.Cells.Item(30,12).Value = 1.0
$ws.Cells.Item(30,13).Value = 49.28281533333333
$ws.Cells.Item(30,14).Value = 147.848446
$ws.Cells.Item(30,15).Value = 0.6169255450395754
$ws.Cells.Item(30,16).Value = 0.6178041890138317
$ws.Cells.Item(30,17).Value = 1.455995068602889
$ws.Cells.Item(30,18).Value = 13.103955617426
$ws.Cells.Item(30,19).Value = 0.00219694852029293
$ws.Cells.Item(30,20).Value = 0.002270135291893619
$ws.Cells.Item(31,5).Value = 1.0
$ws.Cells.Item(31,6).Value = 0.3333333333333333
$ws.Cells.Item(31,7).Value = 0.02954366666666667
$ws.Cells.Item(31,8).Value = 0.088631
$ws.Cells.Item(31,9).Value = 0.003561124252282336
$ws.Cells.Item(31,10).Value = 0.003674522336142324
$ws.Cells.Item(31,13).Value = 2.409433333333333
$ws.Cells.Item(31,14).Value = 7.2283
$ws.Cells.Item(31,15).Value = 0.03016144597968628
$ws.Cells.Item(31,16).Value = 0.03020440282103933
$ws.Cells.Item(31,17).Value = 0.07118349525555555
$ws.Cells.Item(31,18).Value = 0.6406514573000001
$ws.Cells.Item(31,19).Value = 0.0001074086567621643
$ws.Cells.Item(31,20).Value = 0.0001109867528157492

Write-Host "Updated 380 cells"